# Add a new bullet ("An intuitive navigation process to find specific code
# and resources") right after the "A structure for holding organized
# collections of code" bullet on the "Website for Publishing Reusable
# Components" slide, at the same outline level.

$p = $ppt.ActivePresentation

$targetShape = $null
$targetSlide = $null

# Locate the slide / shape that holds the bullet list we need to extend.
foreach ($slide in $p.Slides) {
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -like "*A structure for holding organized collections of code*") {
                $targetSlide = $slide
                $targetShape = $shape
                break
            }
        }
    }
    if ($targetShape -ne $null) { break }
}

$tr = $targetShape.TextFrame.TextRange

# Find which paragraph holds the "A structure for holding organized
# collections of code" bullet, so the new bullet can be inserted right
# after it (paragraph text includes a trailing CR, so trim before
# comparing).
$paraCount = $tr.Paragraphs().Count
$anchorIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $paraText = $tr.Paragraphs($i, 1).Text.TrimEnd([char]13)
    if ($paraText -eq "A structure for holding organized collections of code") {
        $anchorIndex = $i
        break
    }
}

$anchorPara = $tr.Paragraphs($anchorIndex, 1)

# Start a new paragraph right after the anchor bullet (inherits the same
# lvl="1" indent because it sits among the other level-1 bullets).
[void]$anchorPara.InsertAfter([char]13 + "An intuitive navigation process to ")

# Re-fetch the text range / paragraph after the structural edit, then
# append the remaining two runs of the new bullet.
$tr2 = $targetShape.TextFrame.TextRange
$newPara = $tr2.Paragraphs($anchorIndex + 1, 1)
[void]$newPara.InsertAfter("find specific code ")

$tr3 = $targetShape.TextFrame.TextRange
$newPara2 = $tr3.Paragraphs($anchorIndex + 1, 1)
[void]$newPara2.InsertAfter("and resources")
